# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 5 of the
# zh-cn and de-de worksheets, as part of regenerating the handback
# status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-20 08:17:24"
$wsZhCn.Range("G5").Value = "2016-01-20 08:18:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-20 08:17:35"
$wsDeDe.Range("G5").Value = "2016-01-20 08:18:29"
